# Generate Report for Handback
# Updates the handback-status workbook with the newly generated handoff/handback
# file identifiers and timestamps for the two tracked source files.

$wb = $excel.ActiveWorkbook

# Old -> New identifiers
$oldMd1 = "7d23684b-c227-4902-a598-92575ecc296d.md"
$newMd1 = "aa846359-e5b7-4b1a-992e-45eab8e66c07.md"

$oldMd2 = "a865a3d9-612c-4cf4-ba43-b184ebb865f4.md"
$newMd2 = "ffff54f9520a-8ea1-4e33-ba88-650474803040.md"

# New correspond xlf file names (both source files now share the same aa846359 xlf names)
$newXlfZh = "aa846359-e5b7-4b1a-992e-45eab8e66c07.6048720f6a2b7d4f85484e10c32f8fbb6ee781b6.zh-cn.xlf"
$newXlfDe = "aa846359-e5b7-4b1a-992e-45eab8e66c07.6048720f6a2b7d4f85484e10c32f8fbb6ee781b6.de-de.xlf"

# New handoff / handback datetimes
$newHandoffZh = "2016-03-11 22:44:27"
$newHandbackZh = "2016-03-11 22:44:43"

$newHandoffDe = "2016-03-11 22:44:30"
$newHandbackDe = "2016-03-11 22:44:49"

# NOTE: this runtime's PowerShell interpreter does not support *named*
# parameter binding (e.g. "-Sheet $x -Ref $y") for function calls - values
# come through as empty. Positional parameter binding works fine, so all
# helper-function calls below use positional arguments only.
function Set-CellAndLink {
    param($Sheet, $Ref, $NewValue)
    $target = $Sheet.Range($Ref)
    $target.Value = $NewValue
    $targetRow = $target.Row
    $targetCol = $target.Column
    foreach ($l in $Sheet.Hyperlinks) {
        if (($l.Range.Row -eq $targetRow) -and ($l.Range.Column -eq $targetCol)) {
            $l.TextToDisplay = $NewValue
        }
    }
}

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
Set-CellAndLink $wsOverview "A2" $newMd1
Set-CellAndLink $wsOverview "A3" $newMd2

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
# Row 2 (first source file)
Set-CellAndLink $wsZh "A2" $newMd1
Set-CellAndLink $wsZh "D2" $newXlfZh
$wsZh.Range("E2").Value = $newHandoffZh
Set-CellAndLink $wsZh "F2" $newMd1
Set-CellAndLink $wsZh "G2" $newXlfZh
$wsZh.Range("H2").Value = $newHandbackZh
# Row 3 (second source file)
Set-CellAndLink $wsZh "A3" $newMd2
Set-CellAndLink $wsZh "D3" $newXlfZh
$wsZh.Range("E3").Value = $newHandoffZh
Set-CellAndLink $wsZh "F3" $newMd2
Set-CellAndLink $wsZh "G3" $newXlfZh
$wsZh.Range("H3").Value = $newHandbackZh

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
# Row 2 (first source file)
Set-CellAndLink $wsDe "A2" $newMd1
Set-CellAndLink $wsDe "D2" $newXlfDe
$wsDe.Range("E2").Value = $newHandoffDe
Set-CellAndLink $wsDe "F2" $newMd1
Set-CellAndLink $wsDe "G2" $newXlfDe
$wsDe.Range("H2").Value = $newHandbackDe
# Row 3 (second source file)
Set-CellAndLink $wsDe "A3" $newMd2
Set-CellAndLink $wsDe "D3" $newXlfDe
$wsDe.Range("E3").Value = $newHandoffDe
Set-CellAndLink $wsDe "F3" $newMd2
Set-CellAndLink $wsDe "G3" $newXlfDe
$wsDe.Range("H3").Value = $newHandbackDe
